$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift rows 3..12 down to 4..13 (bottom-up copy) to make room for the new row 3
#     (old row 2 "C1..C29 1uF" splits into new row2 "C1" + new row3 "C2" etc.)
#     Row 12 is copied twice (into 12 stays, and into the brand-new row 13) so the
#     extra trailing row inherits real (bordered/filled) formatting instead of blank.
$ws.Range("A12:F12").Copy()
$ws.Range("A13:F13").PasteSpecial(-4104)
$excel.CutCopyMode = $false
for ($r = 12; $r -ge 3; $r--) {
    $src = $ws.Range("A" + $r + ":F" + $r)
    $dst = $ws.Range("A" + ($r+1) + ":F" + ($r+1))
    $src.Copy()
    $dst.PasteSpecial(-4104)
}
$excel.CutCopyMode = $false

# --- Header row: rename LibRef -> MANUF (column E)
$ws.Range("E1").Value = "MANUF"

# Row 2: C1
$ws.Range("A2").Value = "'C1"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "'CAP SMD 0402 1 µF 10V X5R"
$ws.Range("D2").Value = "'CAPC0402N"
$ws.Range("E2").Value = "'WE"
$ws.Range("F2").Value = "'885012105012"
$ws.Rows(2).AutoFit()

# Row 3: C2
$ws.Range("A3").Value = "'C2"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "'CAP SMD 0805 10 µF 10 V X5R"
$ws.Range("D3").Value = "'CAPC0805N"
$ws.Range("E3").Value = "'WE"
$ws.Range("F3").Value = "'885012107010"
$ws.Rows(3).AutoFit()

# Row 4: C3, C4
$ws.Range("A4").Value = "'C3, C4"
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "'CAP SMD 0402 10 pF 10V NP0"
$ws.Range("D4").Value = "'CAPC0402N"
$ws.Range("E4").Value = "'WE"
$ws.Range("F4").Value = "'885012005007"
$ws.Rows(4).AutoFit()

# Row 5: C5, C6, C7, C8, C9, C10, C11, C12, C13, C14, C15, C16, C17, C18, C19, C20, C21, C22, C23, C24, C25, C26, C27, C28, C29
$ws.Range("A5").Value = "'C5, C6, C7, C8, C9, C10, C11, C12, C13, C14, C15, C16, C17, C18, C19, C20, C21, C22, C23, C24, C25, C26, C27, C28, C29"
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = "'CAP SMD 0402 100 nF 10V X5R"
$ws.Range("D5").Value = "'CAPC0402N"
$ws.Range("E5").Value = "'WE"
$ws.Range("F5").Value = "'885012105010"
$ws.Rows(5).RowHeight = 86.4

# Row 6: D1, D2, D3, D4, D5, D6, D7, D8, D9, D10, D11, D12, D13, D15, D16, D17, D18, D19, D20, D21, D22, D23, D24, D25, D26, D27, D28, D29, D30, D31, D32, D33, D34, D35, D36, D37, D38, D39, D40, D41, D42, D43, D44, D45, D46, D47, D48, D49, D50, D51, D52, D53, D54, D55, D56, D57, D58, D59, D60, D61, D62, D63, D64, D65, D66, D67, D68, D69, D70, D71
$ws.Range("A6").Value = "'D1, D2, D3, D4, D5, D6, D7, D8, D9, D10, D11, D12, D13, D15, D16, D17, D18, D19, D20, D21, D22, D23, D24, D25, D26, D27, D28, D29, D30, D31, D32, D33, D34, D35, D36, D37, D38, D39, D40, D41, D42, D43, D44, D45, D46, D47, D48, D49, D50, D51, D52, D53, D54, D55, D56, D57, D58, D59, D60, D61, D62, D63, D64, D65, D66, D67, D68, D69, D70, D71"
$ws.Range("B6").Value = 70
$ws.Range("C6").Value = "'DIODE GEN PURP 80V 125MA SOD323"
$ws.Range("D6").Value = "'SOD323FL"
$ws.Range("E6").Value = "'ON Semi"
$ws.Range("F6").Value = "'1N4148WT"
$ws.Rows(6).RowHeight = 244.8

# Row 7: J1
$ws.Range("A7").Value = "'J1"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = "'USB - mini B USB 2.0 OTG Receptacle Connector 5 Position Surface Mount"
$ws.Range("D7").Value = "'513870530"
$ws.Range("E7").Value = "'MOLEX"
$ws.Range("F7").Value = "'513870530"
$ws.Rows(7).RowHeight = 28.8

# Row 8: LED1, LED2, LED3, LED4, LED5, LED6, LED7, LED8, LED9, LED10, LED11, LED12, LED13, LED14, LED15, LED16, LED17, LED18
$ws.Range("A8").Value = "'LED1, LED2, LED3, LED4, LED5, LED6, LED7, LED8, LED9, LED10, LED11, LED12, LED13, LED14, LED15, LED16, LED17, LED18"
$ws.Range("B8").Value = 18
$ws.Range("C8").Value = "'ADDRESS LED RED/GRN/BLUE"
$ws.Range("D8").Value = "'WS2812"
$ws.Range("E8").Value = "'Worldsemi"
$ws.Range("F8").Value = "'WS2812B"
$ws.Rows(8).RowHeight = 86.4

# Row 9: R1, R2
$ws.Range("A9").Value = "'R1, R2"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "'RES SMD 22 OHM 1% 1/10W 100PPM 0603"
$ws.Range("D9").Value = "'RESC0603N"
$ws.Range("E9").Value = "'VISHAY"
$ws.Range("F9").Value = "'CRCW060322R0FKEA"
$ws.Rows(9).RowHeight = 28.8

# Row 10: R3, R4
$ws.Range("A10").Value = "'R3, R4"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "'RES SMD 4.7K OHM 1% 1/16W 100PPM 0402"
$ws.Range("D10").Value = "'RESC0402N"
$ws.Range("E10").Value = "'VISHAY"
$ws.Range("F10").Value = "'CRCW04024K70FKED"
$ws.Rows(10).RowHeight = 28.8

# Row 11: S72
$ws.Range("A11").Value = "'S72"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "'SWITCH TACTILE SPST-NO 0.05A 12V"
$ws.Range("D11").Value = "'434153017835"
$ws.Range("E11").Value = "'WE"
$ws.Range("F11").Value = "'434153017835"
$ws.Rows(11).AutoFit()

# Row 12: U1
$ws.Range("A12").Value = "'U1"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "'IC MCU 8BIT 32KB FLASH 44VQFN"
$ws.Range("D12").Value = "'44PW_M"
$ws.Range("E12").Value = "'Microchip"
$ws.Range("F12").Value = "'ATMEGA32U4-MU-ND"
$ws.Rows(12).RowHeight = 28.8

# Row 13: Y1
$ws.Range("A13").Value = "'Y1"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "'CRYSTAL 16.0000MHZ 18PF SMD"
$ws.Range("D13").Value = "'WE-XTAL_CFPX-180"
$ws.Range("E13").Value = "'WE/IQD"
$ws.Range("F13").Value = "'830055951"
$ws.Rows(13).AutoFit()

# --- Column widths: C wider (16 -> 35); D single col 19; E:F 16
$ws.Columns("C").ColumnWidth = 34.17
$ws.Columns("D").ColumnWidth = 18.17
$ws.Columns("E:F").ColumnWidth = 15.17

# --- Print scale 69 -> 88 (explicit zoom; also clears the old "fit to page" scaling mode)
$ws.PageSetup.Zoom = 88

